# "matriz SoD" (matrizsod) now loads its header id column ("#") as "X",
# and the sheet becomes the active/selected tab with the cursor resting
# on the first empty row below the data (A8) - as if freshly edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("matrizsod")

# Change header cell A1 from "#" to "X"
$ws.Range("A1").Value = "X"

# Make matrizsod the active sheet and leave the selection on A8
$ws.Activate()
$ws.Range("A8").Select()
